# Auto-generated edit script: updates currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) for specific Leve rows across all 8 job sheets, per the scraped
# market-data refresh described in the commit.

$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 3999
$ws.Range("I46").Value = 3998
$ws.Range("J46").Value = 4000
$ws.Range("K46").Value = 11994
$ws.Range("L46").Value = 12000
$ws.Range("M46").Value = -11875
$ws.Range("N46").Value = -12238

$ws.Range("H60").Value = 3999
$ws.Range("I60").Value = 3998
$ws.Range("J60").Value = 4000
$ws.Range("K60").Value = 11994
$ws.Range("L60").Value = 12000
$ws.Range("M60").Value = -11510
$ws.Range("N60").Value = -12968

$ws.Range("H80").Value = 763.9375
$ws.Range("J80").Value = 962.8889
$ws.Range("L80").Value = 2888.6667
$ws.Range("N80").Value = -4884.6667

$ws.Range("H83").Value = 763.9375
$ws.Range("J83").Value = 962.8889
$ws.Range("L83").Value = 8666.000100000001
$ws.Range("N83").Value = -18650.0001

$ws.Range("H92").Value = 100480.7
$ws.Range("I92").Value = 111478.555
$ws.Range("K92").Value = 111478.555
$ws.Range("M92").Value = -110230.555


# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1725.4193
$ws.Range("I32").Value = 1449.6333
$ws.Range("K32").Value = 1449.6333
$ws.Range("M32").Value = -1162.6333

$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

$ws.Range("H88").Value = 2864.375
$ws.Range("I88").Value = 2407.5
$ws.Range("J88").Value = 3016.6667
$ws.Range("K88").Value = 2407.5
$ws.Range("L88").Value = 3016.6667
$ws.Range("M88").Value = -2001.5
$ws.Range("N88").Value = -3828.6667

$ws.Range("H91").Value = 2864.375
$ws.Range("I91").Value = 2407.5
$ws.Range("J91").Value = 3016.6667
$ws.Range("K91").Value = 2407.5
$ws.Range("L91").Value = 3016.6667
$ws.Range("M91").Value = -1003.5
$ws.Range("N91").Value = -5824.6667

$ws.Range("H122").Value = 4391.143
$ws.Range("I122").Value = 3581.4167
$ws.Range("K122").Value = 10744.2501
$ws.Range("M122").Value = -8294.250100000001

$ws.Range("H132").Value = 2680.375
$ws.Range("I132").Value = 2597.2144
$ws.Range("J132").Value = 3262.5
$ws.Range("K132").Value = 7791.6432
$ws.Range("L132").Value = 9787.5
$ws.Range("M132").Value = -5261.6432
$ws.Range("N132").Value = -14847.5


# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 398.5
$ws.Range("I5").Value = 200
$ws.Range("K5").Value = 200
$ws.Range("M5").Value = -87

$ws.Range("H46").Value = 20000
$ws.Range("I46").Value = 20000
$ws.Range("K46").Value = 20000
$ws.Range("M46").Value = -19702

$ws.Range("H53").Value = 100000
$ws.Range("J53").Value = 100000
$ws.Range("L53").Value = 100000
$ws.Range("N53").Value = -101148


# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2681.6843
$ws.Range("I31").Value = 1335.2
$ws.Range("J31").Value = 7731
$ws.Range("K31").Value = 1335.2
$ws.Range("L31").Value = 7731
$ws.Range("M31").Value = -1040.2
$ws.Range("N31").Value = -8321

$ws.Range("H34").Value = 2681.6843
$ws.Range("I34").Value = 1335.2
$ws.Range("J34").Value = 7731
$ws.Range("K34").Value = 1335.2
$ws.Range("L34").Value = 7731
$ws.Range("M34").Value = -1133.2
$ws.Range("N34").Value = -8135

$ws.Range("H86").Value = 8856
$ws.Range("I86").Value = 8998.25
$ws.Range("J86").Value = 8666.333000000001
$ws.Range("K86").Value = 8998.25
$ws.Range("L86").Value = 8666.333000000001
$ws.Range("M86").Value = -7875.25
$ws.Range("N86").Value = -10912.333

$ws.Range("H89").Value = 8856
$ws.Range("I89").Value = 8998.25
$ws.Range("J89").Value = 8666.333000000001
$ws.Range("K89").Value = 44991.25
$ws.Range("L89").Value = 43331.665
$ws.Range("M89").Value = -39375.25
$ws.Range("N89").Value = -54563.665

$ws.Range("H96").Value = 14767
$ws.Range("J96").Value = 14767
$ws.Range("L96").Value = 14767
$ws.Range("N96").Value = -20259

$ws.Range("H105").Value = 3337.5
$ws.Range("I105").Value = 3113
$ws.Range("J105").Value = 4011
$ws.Range("K105").Value = 3113
$ws.Range("L105").Value = 4011
$ws.Range("M105").Value = -1366
$ws.Range("N105").Value = -7505


# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 1891.0869
$ws.Range("J26").Value = 998
$ws.Range("L26").Value = 2994
$ws.Range("N26").Value = -3570

$ws.Range("H34").Value = 63957.234
$ws.Range("J34").Value = 70458.2
$ws.Range("L34").Value = 211374.6
$ws.Range("N34").Value = -211542.6

$ws.Range("H50").Value = 451.75
$ws.Range("I50").Value = 419
$ws.Range("K50").Value = 1257
$ws.Range("M50").Value = -776

$ws.Range("H53").Value = 451.75
$ws.Range("I53").Value = 419
$ws.Range("K53").Value = 1257
$ws.Range("M53").Value = -776

$ws.Range("H106").Value = 12958.571
$ws.Range("J106").Value = 13535.833
$ws.Range("L106").Value = 40607.499
$ws.Range("N106").Value = -42499.499

$ws.Range("H108").Value = 120.25
$ws.Range("I108").Value = 120.25
$ws.Range("K108").Value = 360.75
$ws.Range("M108").Value = 2519.25

$ws.Range("H137").Value = 4956.6665
$ws.Range("I137").Value = 2946.6667
$ws.Range("J137").Value = 6966.6665
$ws.Range("K137").Value = 8840.000100000001
$ws.Range("L137").Value = 20899.9995
$ws.Range("M137").Value = -3740.000100000001
$ws.Range("N137").Value = -31099.9995

$ws.Range("H140").Value = 771508.5600000001
$ws.Range("I140").Value = 771508.5600000001
$ws.Range("K140").Value = 2314525.68
$ws.Range("M140").Value = -2309345.68


# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 414.35294
$ws.Range("I97").Value = 369.76923
$ws.Range("J97").Value = 559.25
$ws.Range("K97").Value = 369.76923
$ws.Range("L97").Value = 559.25
$ws.Range("M97").Value = 126.23077
$ws.Range("N97").Value = -1551.25

$ws.Range("H107").Value = 2062.3635
$ws.Range("I107").Value = 1798.4286
$ws.Range("J107").Value = 2524.25
$ws.Range("K107").Value = 1798.4286
$ws.Range("L107").Value = 2524.25
$ws.Range("M107").Value = 121.5714
$ws.Range("N107").Value = -6364.25

$ws.Range("H113").Value = 1505.8
$ws.Range("I113").Value = 1382.25
$ws.Range("K113").Value = 1382.25
$ws.Range("M113").Value = 787.75

$ws.Range("H132").Value = 5708.75
$ws.Range("I132").Value = 5037.5
$ws.Range("J132").Value = 7051.25
$ws.Range("K132").Value = 15112.5
$ws.Range("L132").Value = 21153.75
$ws.Range("M132").Value = -12582.5
$ws.Range("N132").Value = -26213.75


# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 26318362
$ws.Range("I136").Value = 2677.4546
$ws.Range("J136").Value = 62502428
$ws.Range("K136").Value = 8032.3638
$ws.Range("L136").Value = 187507284
$ws.Range("M136").Value = -5482.3638
$ws.Range("N136").Value = -187512384


# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 96000
$ws.Range("J124").Value = 96000
$ws.Range("L124").Value = 96000
$ws.Range("N124").Value = -105820

